# Edit script: applies three changes to rela/registro.docx
# 1. Page count "57" -> "67"
# 2. Author name runs (split by spell-check proofErr markers) merged into a single run
# 3. Abstract ("RESUMO") paragraph runs (split by spell-check proofErr markers) merged into a single run
#    (text content is unchanged, only the run/proofErr structure is simplified)

$d = $word.ActiveDocument

# --- Change 1: page count 57 -> 67 -------------------------------------------------
$found1 = $d.Content.Find.Execute("57", $true, $false, $false, $false, $false, $true, 1, $false, "67", 2)

# --- Change 2: author name - merge "Cesar " / "Ryudi" / " " / "Kawakami" runs -------
$found2 = $d.Content.Find.Execute("Cesar Ryudi Kawakami", $true, $false, $false, $false, $false, $true, 1, $false, "Cesar Ryudi Kawakami", 2)

# --- Change 3: abstract paragraph - merge all runs into a single run ---------------
$abstractOld = 'Um sistema online para execução segura de código arbitrário é um sistema computacional voltado para a Internet que permita a execução segura de código proveniente de fontes não-confiáveis. A execução segura de código arbitrário é um tema pouco abordado na área de sistemas operacionais, uma vez que, em geral, assume-se que o usuário será o responsável pelos processos disparados. Os objetivos imediatos deste trabalho são fazer uma exposição dos casos de uso considerados para então elaborar uma arquitetura detalhada de um sistema que resolve o problema e montar um protótipo funcional. São estudadas duas aplicações: competições de programação, que avaliam os competidores executando seus códigos contra uma bateria secreta de testes, e as IDEs online, que permitem a execução de código sem necessitar da instalação de ambientes locais de desenvolvimento. Um levantamento é feito sobre os poucos sistemas já existentes que abordam os mesmos problemas. São descritas as principais tecnologias utilizadas e respectivas fundamentações teóricas durante o trabalho, incluindo Tornado, MongoDB, RabbitMQ, AppArmor e long polling. É apresentado o projeto de um sistema online para execução segura de código arbitrário. São descritas as considerações de design feitas, bem como as estratégias arquiteturais escolhidas, que incluem o uso de programação de alto nível, o uso de escalabilidade horizontal como meio para obtenção de performance e a segurança em profundidade. É apresentada a arquitetura do sistema, com uma exposição sobre a sua visão geral, bem como a arquitetura dos subsistemas, o modelo de segurança e os protocolos de comunicação. Finalmente, um protótipo funcional para cada um dos casos de uso estudados é mostrado, seguido de uma breve análise de validação. O protótipo apresentou-se performante, com escalabilidade horizontal linear e resistente aos ataques testados. No futuro, pode ser feito um aprofundamento sobre as possibilidades criadas por sistemas seguros para execução de código arbitrário, bem como sobre a problemática envolvida na elaboração de sistemas mais completos para o problema estudado.'
$abstractNew = 'Um sistema online para execução segura de código arbitrário é um sistema computacional voltado para a Internet que permita a execução segura de código proveniente de fontes não-confiáveis. A execução segura de código arbitrário é um tema pouco abordado na área de sistemas operacionais, uma vez que, em geral, assume-se que o usuário será o responsável pelos processos disparados. Os objetivos imediatos deste trabalho são fazer uma exposição dos casos de uso considerados para então elaborar uma arquitetura detalhada de um sistema que resolve o problema e montar um protótipo funcional. São estudadas duas aplicações: competições de programação, que avaliam os competidores executando seus códigos contra uma bateria secreta de testes, e as IDEs online, que permitem a execução de código sem necessitar da instalação de ambientes locais de desenvolvimento. Um levantamento é feito sobre os poucos sistemas já existentes que abordam os mesmos problemas. São descritas as principais tecnologias utilizadas e respectivas fundamentações teóricas durante o trabalho, incluindo Tornado, MongoDB, RabbitMQ, AppArmor e long polling. É apresentado o projeto de um sistema online para execução segura de código arbitrário. São descritas as considerações de design feitas, bem como as estratégias arquiteturais escolhidas, que incluem o uso de programação de alto nível, o uso de escalabilidade horizontal como meio para obtenção de performance e a segurança em profundidade. É apresentada a arquitetura do sistema, com uma exposição sobre a sua visão geral, bem como a arquitetura dos subsistemas, o modelo de segurança e os protocolos de comunicação. Finalmente, um protótipo funcional para cada um dos casos de uso estudados é mostrado, seguido de uma breve análise de validação. O protótipo apresentou-se performante, com escalabilidade horizontal linear e resistente aos ataques testados. No futuro, pode ser feito um aprofundamento sobre as possibilidades criadas por sistemas seguros para execução de código arbitrário, bem como sobre a problemática envolvida na elaboração de sistemas mais completos para o problema estudado.'
$found3 = $d.Content.Find.Execute($abstractOld, $true, $false, $false, $false, $false, $true, 1, $false, $abstractNew, 2)

Write-Host "page number replace:" $found1
Write-Host "author name replace:" $found2
Write-Host "abstract replace:" $found3
